$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (dataset shrunk from 5 data rows to 4 data rows)
$ws.Rows("6:6").Delete()

# Adjust column widths (custom accuracy tweak)
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 6.166666666666667
$ws.Columns.Item(22).ColumnWidth = 6.166666666666667

# New dataset values (1000 rows of source data resampled) for rows 2-5
$data = @(
  @(45168.50694444445, 11.794, 7.946, 3.447, 25.671, 18.994, 8.995, 26.449, 14.543, 5.683, 8.23, 10.12, 10.985, 3.013, 9.399, 12.777, 8.54, 2.71, 1.486, 135.674, 25.95, 8.676, 16.478, 8.538, 2.569, 14.652, 7.663, 7.143, 8.202, 10.671, 2.707, 24.031, 4.513, 10.846),
  @(45168.51388888889, 12.328, 8.894, 1.623, 27.096, 21.197, 9.559, 36.499, 15.125, 6.436, 9.285, 10.822, 11.644, 3.14, 9.775, 13.674, 8.624, 1.325, 0.872, 141.458, 27.337, 9.023, 17.924, 9.357, 1.825, 18.364, 7.97, 7.262, 8.485, 11.354, 1.141, 33.609, 4.883, 11.28),
  @(45168.52083333334, 22.449, 16.658, 1.506, 49.091, 39.762, 17.577, 66.759, 27.341, 12.055, 17.745, 19.676, 20.968, 5.676, 17.67, 25.067, 15.018, 1.001, 0.983, 261.643, 49.381, 16.31, 33.052, 17.331, 2.675, 33.005, 14.407, 12.83, 15.064, 20.671, 0.724, 60.839, 9.114, 20.392),
  @(45168.52777777778, 5.66, 4.11, 0.74, 12.53, 9.73, 4.39, 23.33, 6.98, 3.02, 4.22, 5.01, 5.43, 1.46, 4.51, 6.38, 4.03, 0.66, 0.4, 61.4, 12.9, 4.16, 8.48, 4.35, 0.86, 10.82, 3.68, 3.38, 3.95, 5.26, 0.53, 21.78, 2.23, 5.21)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Host "done"